$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 99
$ws.Range("H99").Value = 1655.4286
$ws.Range("J99").Value = 900
$ws.Range("L99").Value = 2700
$ws.Range("N99").Value = -5696
# Row 100
$ws.Range("H100").Value = 2200
$ws.Range("J100").Value = 2200
$ws.Range("L100").Value = 2200
$ws.Range("N100").Value = -3282
# Row 131
$ws.Range("H131").Value = 4239.375
$ws.Range("I131").Value = 986
$ws.Range("K131").Value = 2958
$ws.Range("M131").Value = 2082
# Row 132
$ws.Range("H132").Value = 114158.33
$ws.Range("I132").Value = 274492.38
$ws.Range("K132").Value = 823477.14
$ws.Range("M132").Value = -820947.14
# Row 137
$ws.Range("H137").Value = 5348.136
$ws.Range("I137").Value = 1829.6666
$ws.Range("K137").Value = 5488.9998
$ws.Range("M137").Value = -2938.9998
# Row 138
$ws.Range("H138").Value = 6095.7085
$ws.Range("J138").Value = 8379.406000000001
$ws.Range("L138").Value = 25138.218
$ws.Range("N138").Value = -35418.218

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 19561.21
$ws.Range("I61").Value = 17599.666
$ws.Range("J61").Value = 21326.6
$ws.Range("K61").Value = 17599.666
$ws.Range("L61").Value = 21326.6
$ws.Range("M61").Value = -17387.666
$ws.Range("N61").Value = -21750.6
# Row 74
$ws.Range("H74").Value = 1152.92
$ws.Range("I74").Value = 476
$ws.Range("J74").Value = 1533.6875
$ws.Range("K74").Value = 476
$ws.Range("L74").Value = 1533.6875
$ws.Range("M74").Value = 398
$ws.Range("N74").Value = -3281.6875
# Row 77
$ws.Range("H77").Value = 1152.92
$ws.Range("I77").Value = 476
$ws.Range("J77").Value = 1533.6875
$ws.Range("K77").Value = 2380
$ws.Range("L77").Value = 7668.4375
$ws.Range("M77").Value = 1988
$ws.Range("N77").Value = -16404.4375
# Row 110
$ws.Range("H110").Value = 3844.4443
$ws.Range("I110").Value = 1950
$ws.Range("K110").Value = 1950
$ws.Range("M110").Value = 95
# Row 122
$ws.Range("H122").Value = 4621.8887
$ws.Range("I122").Value = 2749.5386
$ws.Range("K122").Value = 8248.6158
$ws.Range("M122").Value = -5798.6158
# Row 132
$ws.Range("H132").Value = 26591.482
$ws.Range("I132").Value = 34994.824
$ws.Range("K132").Value = 104984.472
$ws.Range("M132").Value = -102454.472
# Row 136
$ws.Range("H136").Value = 19561.21
$ws.Range("I136").Value = 17599.666
$ws.Range("J136").Value = 21326.6
$ws.Range("K136").Value = 52798.99800000001
$ws.Range("L136").Value = 63979.8
$ws.Range("M136").Value = -50248.99800000001
$ws.Range("N136").Value = -69079.79999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 7602.6665
$ws.Range("I20").Value = 6904
$ws.Range("K20").Value = 6904
$ws.Range("M20").Value = -6657
# Row 25
$ws.Range("H25").Value = 20000
$ws.Range("J25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("N25").Value = -20470
# Row 134
$ws.Range("H134").Value = 2236.2642
$ws.Range("I134").Value = 1826.2325
$ws.Range("K134").Value = 5478.6975
$ws.Range("M134").Value = -2943.6975
# Row 138
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Range("H10").Value = 2966.1667
$ws.Range("I10").Value = 1949.25
$ws.Range("K10").Value = 1949.25
$ws.Range("M10").Value = -1810.25
# Row 31
$ws.Range("H31").Value = 2550.0908
$ws.Range("I31").Value = 1144.6154
$ws.Range("J31").Value = 4580.222
$ws.Range("K31").Value = 1144.6154
$ws.Range("L31").Value = 4580.222
$ws.Range("M31").Value = -849.6153999999999
$ws.Range("N31").Value = -5170.222
# Row 34
$ws.Range("H34").Value = 2550.0908
$ws.Range("I34").Value = 1144.6154
$ws.Range("J34").Value = 4580.222
$ws.Range("K34").Value = 1144.6154
$ws.Range("L34").Value = 4580.222
$ws.Range("M34").Value = -942.6153999999999
$ws.Range("N34").Value = -4984.222
# Row 62
$ws.Range("H62").Value = 29879.182
$ws.Range("I62").Value = 2160.8333
$ws.Range("K62").Value = 2160.8333
$ws.Range("M62").Value = -1536.8333
# Row 65
$ws.Range("H65").Value = 29879.182
$ws.Range("I65").Value = 2160.8333
$ws.Range("K65").Value = 10804.1665
$ws.Range("M65").Value = -7684.166499999999
# Row 74
$ws.Range("H74").Value = 79999.5
$ws.Range("J74").Value = 79999.5
$ws.Range("L74").Value = 79999.5
$ws.Range("N74").Value = -81747.5
# Row 77
$ws.Range("H77").Value = 79999.5
$ws.Range("J77").Value = 79999.5
$ws.Range("L77").Value = 239998.5
$ws.Range("N77").Value = -248734.5
# Row 88
$ws.Range("H88").Value = 30681.666
$ws.Range("J88").Value = 30681.666
$ws.Range("L88").Value = 30681.666
$ws.Range("N88").Value = -31493.666
# Row 91
$ws.Range("H91").Value = 30681.666
$ws.Range("J91").Value = 30681.666
$ws.Range("L91").Value = 30681.666
$ws.Range("N91").Value = -33489.666
# Row 99
$ws.Range("H99").Value = 5941.5835
$ws.Range("I99").Value = 2433.3333
$ws.Range("J99").Value = 7111
$ws.Range("K99").Value = 2433.3333
$ws.Range("L99").Value = 7111
$ws.Range("M99").Value = -935.3332999999998
$ws.Range("N99").Value = -10107
# Row 126
$ws.Range("H126").Value = 5941.5835
$ws.Range("I126").Value = 2433.3333
$ws.Range("J126").Value = 7111
$ws.Range("K126").Value = 7299.999899999999
$ws.Range("L126").Value = 21333
$ws.Range("M126").Value = -4829.999899999999
$ws.Range("N126").Value = -26273
# Row 132
$ws.Range("H132").Value = 9265835
$ws.Range("I132").Value = 10756318
$ws.Range("K132").Value = 32268954
$ws.Range("M132").Value = -32266424
# Row 134
$ws.Range("H134").Value = 2165.0244
$ws.Range("I134").Value = 2165.0244
$ws.Range("K134").Value = 6495.073199999999
$ws.Range("M134").Value = -3960.073199999999

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 10997.889
$ws.Range("I3").Value = 3996.8333
$ws.Range("J3").Value = 25000
$ws.Range("K3").Value = 11990.4999
$ws.Range("L3").Value = 75000
$ws.Range("M3").Value = -11878.4999
$ws.Range("N3").Value = -75224
# Row 122
$ws.Range("H122").Value = 546
$ws.Range("I122").Value = 503.66666
$ws.Range("J122").Value = 564.1429000000001
$ws.Range("K122").Value = 4532.99994
$ws.Range("L122").Value = 5077.2861
$ws.Range("M122").Value = -2082.99994
$ws.Range("N122").Value = -9977.286100000001
# Row 132
$ws.Range("H132").Value = 2270.9092
$ws.Range("I132").Value = 1682.2
$ws.Range("K132").Value = 15139.8
$ws.Range("M132").Value = -12609.8
# Row 141
$ws.Range("H141").Value = 12422.8
$ws.Range("I141").Value = 5861.4546
$ws.Range("J141").Value = 20442.223
$ws.Range("K141").Value = 17584.3638
$ws.Range("L141").Value = 61326.66900000001
$ws.Range("M141").Value = -12404.3638
$ws.Range("N141").Value = -71686.66900000001

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 6181.077
$ws.Range("I132").Value = 5994
$ws.Range("J132").Value = 6804.6665
$ws.Range("K132").Value = 17982
$ws.Range("L132").Value = 20413.9995
$ws.Range("M132").Value = -15452
$ws.Range("N132").Value = -25473.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1449.3334
$ws.Range("I22").Value = 583
$ws.Range("J22").Value = 2315.6667
$ws.Range("K22").Value = 583
$ws.Range("L22").Value = 2315.6667
$ws.Range("M22").Value = -288
$ws.Range("N22").Value = -2905.6667
# Row 27
$ws.Range("H27").Value = 1449.3334
$ws.Range("I27").Value = 583
$ws.Range("J27").Value = 2315.6667
$ws.Range("K27").Value = 583
$ws.Range("L27").Value = 2315.6667
$ws.Range("M27").Value = -476
$ws.Range("N27").Value = -2529.6667
# Row 40
$ws.Range("H40").Value = 5174.1763
$ws.Range("I40").Value = 4095
$ws.Range("J40").Value = 7152.6665
$ws.Range("K40").Value = 4095
$ws.Range("L40").Value = 7152.6665
$ws.Range("M40").Value = -3959
$ws.Range("N40").Value = -7424.6665
# Row 122
$ws.Range("H122").Value = 51952360
$ws.Range("I122").Value = 66670120
$ws.Range("J122").Value = 20414302
$ws.Range("K122").Value = 200010360
$ws.Range("L122").Value = 61242906
$ws.Range("M122").Value = -200007910
$ws.Range("N122").Value = -61247806
# Row 132
$ws.Range("H132").Value = 3705.9604
$ws.Range("I132").Value = 2817.4312
$ws.Range("J132").Value = 6569
$ws.Range("K132").Value = 8452.293600000001
$ws.Range("L132").Value = 19707
$ws.Range("M132").Value = -5922.293600000001
$ws.Range("N132").Value = -24767
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5642.3335
$ws.Range("I62").Value = 3357.4
$ws.Range("K62").Value = 3357.4
$ws.Range("M62").Value = -2733.4
# Row 65
$ws.Range("H65").Value = 5642.3335
$ws.Range("I65").Value = 3357.4
$ws.Range("K65").Value = 16787
$ws.Range("M65").Value = -13667
# Row 122
$ws.Range("H122").Value = 2649.7627
$ws.Range("I122").Value = 1980.3265
$ws.Range("J122").Value = 5930
$ws.Range("K122").Value = 5940.979499999999
$ws.Range("L122").Value = 17790
$ws.Range("M122").Value = -3490.979499999999
$ws.Range("N122").Value = -22690
# Row 132
$ws.Range("H132").Value = 17245686
$ws.Range("I132").Value = 1796.2916
$ws.Range("K132").Value = 5388.8748
$ws.Range("M132").Value = -2858.8748
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 138
$ws.Range("H138").Value = 79974.5
$ws.Range("J138").Value = 79974.5
$ws.Range("L138").Value = 79974.5
$ws.Range("N138").Value = -90254.5

